$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '90.591.31'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.56%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.135.76'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.11%  '

$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.06'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.46%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '622.25'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.52%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.13'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +27.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.364'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.31%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.01%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.134.69'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.19%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.738'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.74%  '

$ws.Range("E12").Value = '  +6.00%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000246'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.66%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.66'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.83%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '35.24'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +6.42%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '90.285.93'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.84%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.712.43'
$ws.Range("D17").Style = "Normal"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.138.81'
$ws.Range("D18").Style = "Normal"

$ws.Range("E19").Value = '  +4.13%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.56'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.21%  '

$ws.Range("E21").Value = '  -9.81%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '463.05'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +7.35%  '

$ws.Range("B23").Value = 'Polkadot'
$ws.Range("C23").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.78'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +12.25%  '

$ws.Range("B24").Value = 'Uniswap'
$ws.Range("C24").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.07'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.96%  '

$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '95.08'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +13.54%  '

$ws.Range("B26").Value = 'NEARProtocol'
$ws.Range("C26").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '5.78'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.40%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.30'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.22%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.307.36'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.36%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("D29").Style = "Normal"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.164'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.51%  '

$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.220'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +51.10%  '

$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '9.24'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +6.30%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.75'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +16.22%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '518.40'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.05%  '

$ws.Range("E35").Value = '  +3.96%  '

$ws.Range("E36").Value = '  +5.20%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.02'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.87%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.33'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.28%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.60'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -7.72%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0919'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +28.06%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.428'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +15.25%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '22.23'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.31%  '

$ws.Range("E43").Value = '  -24.80%  '

$ws.Range("E44").Value = '  +0.06%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.99'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +6.00%  '

$ws.Range("E46").Value = '  +0.00%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.728'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +19.77%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.76'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +12.75%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '150.44'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.73%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.37'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +8.22%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '45.29'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.48%  '
